$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: tensorflow.blog entry title/link update
$ws.Range("D12").Value = "“코딩 뇌를 깨우는 파이썬”이 출간되었습니다!"
$ws.Range("E12").Value = "https://tensorflow.blog/2023/08/24/%ec%bd%94%eb%94%a9-%eb%87%8c%eb%a5%bc-%ea%b9%a8%ec%9a%b0%eb%8a%94-%ed%8c%8c%ec%9d%b4%ec%8d%ac%ec%9d%b4-%ec%b6%9c%ea%b0%84%eb%90%98%ec%97%88%ec%8a%b5%eb%8b%88%eb%8b%a4/"

# Row 36: dmqm_seminar entry title/link update
$ws.Range("D36").Value = "Image Editing with Diffusion Model"
$ws.Range("E36").Value = "http://dmqm.korea.ac.kr/activity/seminar/420"

# Row 51: bskyvsion entry title/link update
$ws.Range("D51").Value = "[python] 딕셔너리를 json으로 변환하는 방법"
$ws.Range("E51").Value = "https://bskyvision.com/entry/python-%EB%94%95%EC%85%94%EB%84%88%EB%A6%AC%EB%A5%BC-json%EC%9C%BC%EB%A1%9C-%EB%B3%80%ED%99%98%ED%95%98%EB%8A%94-%EB%B0%A9%EB%B2%95"
